# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.281.08"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "2.285.92"
$ws.Range("E3").Value = "  -0.79%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0902"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.12%  "

$ws.Range("E14").Value = "  -0.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.44%  "

$ws.Range("D16").Value = "2.632.35"
$ws.Range("E16").Value = "  -0.82%  "

$ws.Range("D17").Value = "2.287.55"
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("D18").Value = "42.368.84"
$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("E19").Value = "  -5.49%  "

$ws.Range("E20").Value = "  -0.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +30.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.84%  "

$ws.Range("E26").Value = "  +0.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.91%  "

$ws.Range("E33").Value = "  -1.02%  "

$ws.Range("E34").Value = "  +0.34%  "

$ws.Range("E35").Value = "  -2.43%  "

$ws.Range("E36").Value = "  -14.20%  "

$ws.Range("E37").Value = "  -1.16%  "

$ws.Range("E38").Value = "  +0.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.40%  "

$ws.Range("E40").Value = "  -7.23%  "

$ws.Range("E41").Value = "  +1.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.59%  "

$ws.Range("E43").Value = "  +0.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.224"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "80.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.54%  "

$ws.Range("D51").Value = "1.597.12"
$ws.Range("E51").Value = "  +2.91%  "
